# Refresh the cryptos list: updated Price (D) / Volume(1h) (E) figures,
# and WEMIXToken (row 33) now ranks above Filecoin (row 34).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range('D2').Value = '42.960.10'
$ws.Range('E2').Value = '  -0.48%  '
# Row 3 - Ethereum
$ws.Range('D3').Value = '2.369.96'
$ws.Range('E3').Value = '  -1.43%  '
# Row 5 - BNB
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '318.14'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -5.06%  '
# Row 6 - Solana
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '107.78'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +2.00%  '
# Row 7 - XRP
$ws.Range('E7').Value = '  -2.31%  '
# Row 8 - USDC
$ws.Range('E8').Value = '  +0.00%  '
# Row 9 - Cardano
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.623'
$ws.Range('D9').ClearFormats()
# Row 10 - Avalanche
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '41.91'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -1.08%  '
# Row 11 - Dogecoin
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0931'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -0.97%  '
# Row 12 - Polkadot
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '8.58'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -1.96%  '
# Row 13 - Polygon
$ws.Range('E13').Value = '  -4.67%  '
# Row 14 - TRON
$ws.Range('E14').Value = '  +0.06%  '
# Row 15 - Chainlink
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '16.19'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -6.54%  '
# Row 16 - WrappedliquidstakedEther2.0
$ws.Range('D16').Value = '2.729.75'
$ws.Range('E16').Value = '  -1.33%  '
# Row 17 - WrappedEther
$ws.Range('D17').Value = '2.378.99'
$ws.Range('E17').Value = '  -0.87%  '
# Row 18 - WrappedBTC
$ws.Range('D18').Value = '42.913.89'
$ws.Range('E18').Value = '  -0.66%  '
# Row 19 - Uniswap
$ws.Range('E19').Value = '  -0.43%  '
# Row 20 - ShibaInu
$ws.Range('E20').Value = '  -2.18%  '
# Row 21 - Litecoin
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '76.24'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -1.26%  '
# Row 22 - PancakeSwap
$ws.Range('E22').Value = '  -5.07%  '
# Row 23 - BitcoinCash
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '257.70'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -6.66%  '
# Row 24 - ImmutableX
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.34'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -3.33%  '
# Row 25 - InternetComputer(DFINITY)
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.42'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -4.80%  '
# Row 26 - Dai
$ws.Range('E26').Value = '  +0.13%  '
# Row 27 - Cosmos
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.51'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -3.77%  '
# Row 28 - EthereumClassic
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '23.06'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -0.92%  '
# Row 29 - Toncoin
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.25'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +2.51%  '
# Row 30 - InjectiveProtocol
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '37.08'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -0.26%  '
# Row 31 - Monero
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '171.13'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -2.44%  '
# Row 32 - Hedera
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.0900'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -4.38%  '
# Row 33 - WEMIXToken (was Filecoin - rows 33/34 swapped)
$ws.Range('B33').Value = 'WEMIXToken'
$ws.Range('C33').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.95'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -6.57%  '
# Row 34 - Filecoin (was WEMIXToken - rows 33/34 swapped)
$ws.Range('B34').Value = 'Filecoin'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.05'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -0.95%  '
# Row 35 - Kaspa
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.123'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +12.64%  '
# Row 36 - Stellar
$ws.Range('E36').Value = '  -3.22%  '
# Row 37 - RenderToken
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.75'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -3.03%  '
# Row 38 - VeChain
$ws.Range('E38').Value = '  -0.29%  '
# Row 39 - NEARProtocol
$ws.Range('E39').Value = '  -4.50%  '
# Row 40 - LidoDAOToken
$ws.Range('E40').Value = '  -4.47%  '
# Row 41 - ARBITRUM
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.54'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -1.48%  '
# Row 42 - Algorand
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.242'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +2.34%  '
# Row 43 - MultiversX
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '71.63'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +1.47%  '
# Row 44 - FirstDigitalUSD
$ws.Range('E44').Value = '  -0.05%  '
# Row 45 - Celestia
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '12.40'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +0.34%  '
# Row 46 - BitcoinSV
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '88.91'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -3.98%  '
# Row 47 - Aave
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '113.23'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -7.09%  '
# Row 48 - THORChain
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '5.57'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -0.66%  '
# Row 49 - FraxShare
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.25'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -0.21%  '
# Row 50 - ordi
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '77.41'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +7.91%  '
# Row 51 - TrustWalletToken
$ws.Range('E51').Value = '  -1.78%  '
